$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 22:06"

# Swap Siria / Surinam rows (row 133 becomes Siria with fresh data,
# row 134 becomes Surinam carrying the data that used to belong to Surinam)
$ws.Range("A133").Value = "Siria"
$ws.Range("A134").Value = "Surinam"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8428377
$ws.Range("C4").Value = 36650
$ws.Range("D4").Value = 5487758
$ws.Range("E4").Value = 2715575
$ws.Range("G4").Value = 300
$ws.Range("H4").Value = 225044

# Row 21 - Alemania
$ws.Range("B21").Value = 373821
$ws.Range("C21").Value = 6840
$ws.Range("D21").Value = 294800
$ws.Range("E21").Value = 69125
$ws.Range("G21").Value = 30
$ws.Range("H21").Value = 9896

# Row 24 - Turquia
$ws.Range("B24").Value = 349519
$ws.Range("C24").Value = 2026
$ws.Range("E24").Value = 34721

# Row 117
$ws.Range("B117").Value = 7800
$ws.Range("C117").Value = 48
$ws.Range("D117").Value = 6620
$ws.Range("E117").Value = 1093
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 87

# Row 133 - Siria (new data)
$ws.Range("B133").Value = 5134
$ws.Range("C133").Value = 57
$ws.Range("D133").Value = 1565
$ws.Range("E133").Value = 3318
$ws.Range("G133").Value = 3
$ws.Range("H133").Value = 251

# Row 134 - Surinam (carries former Surinam data)
$ws.Range("B134").Value = 5130
$ws.Range("D134").Value = 4944
$ws.Range("E134").Value = 77
$ws.Range("H134").Value = 109

# Row 140
$ws.Range("B140").Value = 4334
$ws.Range("C140").Value = 12
$ws.Range("D140").Value = 4040
$ws.Range("E140").Value = 260
